# This script applies the text corrections found in the diff to sheet "upiti"
# column A, and updates the active selection to A100 (matching the final
# sheetView state in the target workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("upiti")

# Map of row -> corrected text (typo / wording fixes applied in the edit)
$edits = @{
    73  = "ekstrakcija podataka iz html sadrzaja"
    77  = "podvuci tekst u nazivu vidzeta"
    79  = "kopiranje fajla u odredjenu putanju"
    87  = "konverzija uin8 niza u sliku"
    90  = "kako citati sadrzaj .gz kompresovanog fajla"
    92  = "ekstrakcija podataka iz tekstualnog fajla"
    93  = "pozicija pod stringova u stringu"
    94  = "citanje elemenata from html-a <td>"
    97  = "parsiraj string upit u urla"
    100 = "kako citati .csv fajl na efikasan nacin"
}

foreach ($row in $edits.Keys) {
    $ws.Cells.Item($row, 1).Value = $edits[$row]
}

# Update the view: select A100 (this also naturally scrolls the window so
# that A100 is visible, matching the removal of the explicit topLeftCell).
$ws.Range("A100").Select()
